$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.182.59'
$ws.Range('E2').Value = '  -1.31%  '
$ws.Range('D3').Value = '2.173.48'
$ws.Range('E3').Value = '  -2.30%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '250.41'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('E6').Value = '  -2.91%  '
$ws.Range('E7').Value = '  -7.98%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.576'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -3.20%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '58.83'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.88%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '36.22'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -11.73%  '
$ws.Range('E12').Value = '  -3.54%  '
$ws.Range('E13').Value = '  -1.57%  '
$ws.Range('E14').Value = '  -5.18%  '
$ws.Range('D15').Value = '2.498.04'
$ws.Range('E15').Value = '  -2.19%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.25'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -4.77%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.844'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.37%  '
$ws.Range('D18').Value = '2.172.53'
$ws.Range('E18').Value = '  -2.40%  '
$ws.Range('D19').Value = '41.089.51'
$ws.Range('E19').Value = '  -1.54%  '
$ws.Range('D20').Value = '0.0₃0944'
$ws.Range('E20').Value = '  -2.28%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '71.53'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.86%  '
$ws.Range('E22').Value = '  -2.87%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '229.67'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.41%  '
$ws.Range('E24').Value = '  -4.83%  '
$ws.Range('E25').Value = '  -5.94%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.41'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +6.53%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('E28').Value = '  -5.14%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.11'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '168.34'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.65%  '
$ws.Range('E32').Value = '  -2.72%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.69'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.15%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0750'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.29%  '
$ws.Range('E35').Value = '  -3.43%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.49'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.85%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.92'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.99%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '24.33'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -5.05%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0305'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.83%  '
$ws.Range('B40').Value = 'FTXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.43'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +11.61%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.20'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.85%  '
$ws.Range('E42').Value = '  -7.90%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '11.34'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -6.29%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '60.63'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -9.09%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.45'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.71%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0990'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.02%  '
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('E48').Value = '  -6.99%  '
$ws.Range('E49').Value = '  -2.51%  '
$ws.Range('E50').Value = '  -4.64%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.18'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -10.36%  '
